# Fixing the big mistake: correct the "Total" (column B) and "Community"
# (column D) monthly consumption values on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New corrected values for column B (Total), rows 2-13
$totalValues = @(
    7272.422255400007,
    6760.700009166674,
    7236.999968266674,
    7006.445791766673,
    7266.316580000007,
    6990.85414203334,
    7238.452864566674,
    7266.884610816674,
    6977.038504816675,
    7246.483091233341,
    6999.71887263334,
    6964.543750066674
)

# New corrected values for column D (Community), rows 2-13
$communityValues = @(
    353.5872912166667,
    328.8478758333333,
    349.7912609333333,
    331.3696418,
    353.7159776,
    334.4684250833333,
    347.1777969333334,
    355.2152107833333,
    344.4895291833333,
    350.16805895,
    341.6009272833333,
    329.9399408166667
)

for ($i = 0; $i -lt $totalValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $totalValues[$i]
    $ws.Cells.Item($row, 4).Value = $communityValues[$i]
}

$wb.Save()
